$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-11 from 2023-10-09 (45208) to 2023-10-13 (45212)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
